$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (Changed) date column C for rows 2-5 from 45221 to 45224
$ws.Range("C2:C5").Value = 45224
